$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2449
$ws.Range("F3").Value = 718
$ws.Range("F4").Value = 237
$ws.Range("F5").Value = 406
$ws.Range("F6").Value = 679
$ws.Range("F8").Value = 882
$ws.Range("F9").Value = 560
$ws.Range("F10").Value = 923
$ws.Range("F14").Value = 47
$ws.Range("F16").Value = 1061
$ws.Range("F17").Value = 23828
$ws.Range("G17").Value = "暂时售罄"
$ws.Range("F18").Value = 2175
$ws.Range("F19").Value = 139
$ws.Range("F20").Value = 351
$ws.Range("F21").Value = 26
$ws.Range("F22").Value = 46
$ws.Range("F23").Value = 347
$ws.Range("F24").Value = 205
$ws.Range("F25").Value = 62
$ws.Range("F26").Value = 226
$ws.Range("F28").Value = 47
$ws.Range("F30").Value = 335
$ws.Range("F32").Value = 430
$ws.Range("F33").Value = 185

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 187
$ws.Range("F7").Value = 250
$ws.Range("F8").Value = 17
$ws.Range("F10").Value = 3588
$ws.Range("F12").Value = 141
$ws.Range("F16").Value = 14
$ws.Range("F17").Value = 133

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 279
$ws.Range("F3").Value = 157
$ws.Range("F4").Value = 741

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 279
$ws.Range("F3").Value = 157
$ws.Range("F4").Value = 2449
$ws.Range("F5").Value = 741
$ws.Range("F6").Value = 718
$ws.Range("F7").Value = 237
$ws.Range("F8").Value = 407
$ws.Range("F9").Value = 679
$ws.Range("F11").Value = 187
$ws.Range("F14").Value = 250
$ws.Range("F16").Value = 882
$ws.Range("F17").Value = 560
$ws.Range("F18").Value = 923
$ws.Range("F21").Value = 47
$ws.Range("F23").Value = 1061
$ws.Range("F24").Value = 23829
$ws.Range("G24").Value = "暂时售罄"
$ws.Range("F25").Value = 17
$ws.Range("F28").Value = 141
$ws.Range("F30").Value = 2175
$ws.Range("F31").Value = 139
$ws.Range("F32").Value = 351
$ws.Range("F33").Value = 26
$ws.Range("F36").Value = 347
$ws.Range("F37").Value = 205
$ws.Range("F38").Value = 62
$ws.Range("F39").Value = 226
$ws.Range("F41").Value = 14
$ws.Range("F42").Value = 47
$ws.Range("F44").Value = 133
$ws.Range("F47").Value = 430
$ws.Range("F48").Value = 185
